# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary fields -------------------------------------------------
# Valor Mora total (E11): 47 old periods dropped, 48 new periods of the same
# worker/value (36341 each) -> 48 * 36341 = 1744368
$ws.Range("E11").Value = 1744368

# Cant. Trabajadores (C13): only one worker remains now
$ws.Range("C13").Value = 1

# Cant. Periodos (F13): 48 periods now (2108 .. 2507)
$ws.Range("F13").Value = 48

# --- Detail rows 16-63 (period list), newest (2507) first, descending -----
$periods = @(
    "2507","2506","2505","2504","2503","2502","2501",
    "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401",
    "2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301",
    "2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201",
    "2112","2111","2110","2109","2108"
)

$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 5).Value = $p
    $row = $row + 1
}

# --- Row 63 now reuses the single remaining worker (previously the last
# table row belonged to a different worker/debt that has been removed) ----
$ws.Range("C63").Value = "1052088371"
$ws.Range("D63").Value = "EVA SANDRITH GARCIA TERAN"
$ws.Range("E63").Value = "2108"
$ws.Range("F63").Value = 36341
$ws.Range("G63").Value = 908526

# --- Column D width shrinks now that the longer worker name is gone -------
$ws.Columns.Item(4).ColumnWidth = 27.1
